# Finish Service,Testimonial by Tham Vinh Thanh
# Fill in the Build start/finish dates and % Build completion for the
# three tasks that were just completed: "Latest Service" (row 13),
# "Clients Testimonial" (row 17) and "Get Connected" (row 19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest Service - finished 2024-10-23
$ws.Range("B13").Value = 45588
$ws.Range("C13").Value = 45588
$ws.Range("D13").Value = 1

# Clients Testimonial - finished 2024-10-24
$ws.Range("B17").Value = 45589
$ws.Range("C17").Value = 45589
$ws.Range("D17").Value = 1

# Get Connected - finished 2024-10-24
$ws.Range("B19").Value = 45589
$ws.Range("C19").Value = 45589
$ws.Range("D19").Value = 1

# Leave the selection where the author left off editing
$ws.Range("G19").Select()
